$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 24550
$ws.Range("I62").Value = 17906
$ws.Range("J62").Value = 55555.332
$ws.Range("K62").Value = 17906
$ws.Range("L62").Value = 55555.332
$ws.Range("M62").Value = -17282
$ws.Range("N62").Value = -56803.332

$ws.Range("H65").Value = 24550
$ws.Range("I65").Value = 17906
$ws.Range("J65").Value = 55555.332
$ws.Range("K65").Value = 89530
$ws.Range("L65").Value = 277776.66
$ws.Range("M65").Value = -86410
$ws.Range("N65").Value = -284016.66

$ws.Range("H75").Value = 0
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").ClearContents()

$ws.Range("H76").Value = 4852.75
$ws.Range("I76").Value = 4999

$ws.Range("H78").Value = 0
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("K78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").ClearContents()

$ws.Range("H79").Value = 4852.75
$ws.Range("I79").Value = 4999

$ws.Range("I99").Value = 1499
$ws.Range("J99").Value = 1955.8572
$ws.Range("K99").Value = 4497
$ws.Range("L99").Value = 5867.571599999999
$ws.Range("M99").Value = -2999
$ws.Range("N99").Value = -8863.571599999999

$ws.Range("H112").Value = 1723.9048
$ws.Range("I112").Value = 1006.2857
$ws.Range("J112").Value = 2082.7144
$ws.Range("K112").Value = 3018.8571
$ws.Range("L112").Value = 6248.1432
$ws.Range("M112").Value = -1910.8571
$ws.Range("N112").Value = -8464.143199999999

$ws.Range("H124").Value = 70000
$ws.Range("I124").Value = 0
$ws.Range("J124").Value = 70000
$ws.Range("K124").Value = 0
$ws.Range("L124").Value = 70000
$ws.Range("N124").Value = -79820

$ws.Range("H127").Value = 1917
$ws.Range("I127").Value = 1329
$ws.Range("J127").Value = 2799
$ws.Range("K127").Value = 3987
$ws.Range("L127").Value = 8397
$ws.Range("M127").Value = 973
$ws.Range("N127").Value = -18317

$ws.Range("H129").Value = 2396.7144
$ws.Range("I129").Value = 2244.625
$ws.Range("J129").Value = 2883.4
$ws.Range("K129").Value = 6733.875
$ws.Range("L129").Value = 8650.200000000001
$ws.Range("M129").Value = -1733.875
$ws.Range("N129").Value = -18650.2

$ws.Range("H132").Value = 501293.44
$ws.Range("I132").Value = 563283.6
$ws.Range("J132").Value = 5372.25
$ws.Range("K132").Value = 1689850.8
$ws.Range("L132").Value = 16116.75
$ws.Range("M132").Value = -1687320.8
$ws.Range("N132").Value = -21176.75

$ws.Range("H137").Value = 2302.318
$ws.Range("I137").Value = 1570.5883
$ws.Range("J137").Value = 2763.037
$ws.Range("K137").Value = 4711.7649
$ws.Range("L137").Value = 8289.110999999999
$ws.Range("M137").Value = -2161.7649
$ws.Range("N137").Value = -13389.111

$ws.Range("H138").Value = 3391.3635
$ws.Range("I138").Value = 2485.2
$ws.Range("J138").Value = 3785.348
$ws.Range("K138").Value = 7455.599999999999
$ws.Range("L138").Value = 11356.044
$ws.Range("M138").Value = -2315.599999999999
$ws.Range("N138").Value = -21636.044

$ws.Range("H141").Value = 1087.5333
$ws.Range("I141").Value = 1022.3571
$ws.Range("J141").Value = 2000
$ws.Range("K141").Value = 3067.0713
$ws.Range("L141").Value = 6000
$ws.Range("M141").Value = 2112.9287
$ws.Range("N141").Value = -16360

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6953917
$ws.Range("I32").Value = 9809799
$ws.Range("J32").Value = 18203.38
$ws.Range("K32").Value = 9809799
$ws.Range("L32").Value = 18203.38
$ws.Range("M32").Value = -9809512
$ws.Range("N32").Value = -18777.38

$ws.Range("H45").Value = 3674.1765
$ws.Range("I45").Value = 3490.8
$ws.Range("J45").Value = 5049.5
$ws.Range("K45").Value = 3490.8
$ws.Range("L45").Value = 5049.5
$ws.Range("M45").Value = -3113.8
$ws.Range("N45").Value = -5803.5

$ws.Range("H50").Value = 1637.7
$ws.Range("I50").Value = 2866
$ws.Range("J50").Value = 1111.2858
$ws.Range("K50").Value = 2866
$ws.Range("L50").Value = 1111.2858
$ws.Range("M50").Value = -2152
$ws.Range("N50").Value = -2539.2858

$ws.Range("H61").Value = 8508.208000000001
$ws.Range("I61").Value = 4313.467
$ws.Range("J61").Value = 15499.444
$ws.Range("K61").Value = 4313.467
$ws.Range("L61").Value = 15499.444
$ws.Range("M61").Value = -4101.467
$ws.Range("N61").Value = -15923.444

$ws.Range("H132").Value = 974281.0600000001
$ws.Range("I132").Value = 1111464
$ws.Range("J132").Value = 14000
$ws.Range("K132").Value = 3334392
$ws.Range("L132").Value = 42000
$ws.Range("M132").Value = -3331862
$ws.Range("N132").Value = -47060

$ws.Range("H136").Value = 8508.208000000001
$ws.Range("I136").Value = 4313.467
$ws.Range("J136").Value = 15499.444
$ws.Range("K136").Value = 12940.401
$ws.Range("L136").Value = 46498.33199999999
$ws.Range("M136").Value = -10390.401
$ws.Range("N136").Value = -51598.33199999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2152.45
$ws.Range("I20").Value = 2083.8
$ws.Range("J20").Value = 2358.4
$ws.Range("K20").Value = 2083.8
$ws.Range("L20").Value = 2358.4
$ws.Range("M20").Value = -1836.8
$ws.Range("N20").Value = -2852.4

$ws.Range("H35").Value = 93798
$ws.Range("I35").Value = 0
$ws.Range("J35").Value = 93798
$ws.Range("K35").Value = 0
$ws.Range("L35").Value = 93798
$ws.Range("N35").Value = -94418

$ws.Range("H80").Value = 6174.385
$ws.Range("I80").Value = 13228.5
$ws.Range("J80").Value = 3039.2222
$ws.Range("K80").Value = 13228.5
$ws.Range("L80").Value = 3039.2222
$ws.Range("M80").Value = -12230.5
$ws.Range("N80").Value = -5035.2222

$ws.Range("H83").Value = 6174.385
$ws.Range("I83").Value = 13228.5
$ws.Range("J83").Value = 3039.2222
$ws.Range("K83").Value = 66142.5
$ws.Range("L83").Value = 15196.111
$ws.Range("M83").Value = -61150.5
$ws.Range("N83").Value = -25180.111

$ws.Range("H88").Value = 49959.145
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 49959.145
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 49959.145
$ws.Range("N88").Value = -50771.145

$ws.Range("H91").Value = 49959.145
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 49959.145
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 49959.145
$ws.Range("N91").Value = -52767.145

$ws.Range("H95").Value = 38187
$ws.Range("I95").Value = 0
$ws.Range("J95").Value = 38187
$ws.Range("K95").Value = 0
$ws.Range("L95").Value = 38187
$ws.Range("N95").Value = -43679

$ws.Range("H107").Value = 3744.7932
$ws.Range("I107").Value = 3253.8076
$ws.Range("J107").Value = 8000
$ws.Range("K107").Value = 3253.8076
$ws.Range("L107").Value = 8000
$ws.Range("M107").Value = -1333.8076
$ws.Range("N107").Value = -11840

$ws.Range("H134").Value = 923886.5
$ws.Range("I134").Value = 1185287.8
$ws.Range("J134").Value = 8982
$ws.Range("K134").Value = 3555863.4
$ws.Range("L134").Value = 26946
$ws.Range("M134").Value = -3553328.4
$ws.Range("N134").Value = -32016

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H25").Value = 9415.75
$ws.Range("I25").Value = 4221
$ws.Range("J25").Value = 25000
$ws.Range("K25").Value = 4221
$ws.Range("L25").Value = 25000
$ws.Range("M25").Value = -4047
$ws.Range("N25").Value = -25348

$ws.Range("H31").Value = 7981.48
$ws.Range("I31").Value = 1592
$ws.Range("J31").Value = 10466.277
$ws.Range("K31").Value = 1592
$ws.Range("L31").Value = 10466.277
$ws.Range("M31").Value = -1297
$ws.Range("N31").Value = -11056.277

$ws.Range("H34").Value = 7981.48
$ws.Range("I34").Value = 1592
$ws.Range("J34").Value = 10466.277
$ws.Range("K34").Value = 1592
$ws.Range("L34").Value = 10466.277
$ws.Range("M34").Value = -1390
$ws.Range("N34").Value = -10870.277

$ws.Range("H122").Value = 2868.5
$ws.Range("I122").Value = 2052.75
$ws.Range("J122").Value = 4500
$ws.Range("K122").Value = 6158.25
$ws.Range("L122").Value = 13500
$ws.Range("M122").Value = -3708.25
$ws.Range("N122").Value = -18400

$ws.Range("H132").Value = 5565903.5
$ws.Range("I132").Value = 11517.942
$ws.Range("J132").Value = 25006252
$ws.Range("K132").Value = 34553.826
$ws.Range("L132").Value = 75018756
$ws.Range("M132").Value = -32023.826
$ws.Range("N132").Value = -75023816

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 8394.939
$ws.Range("I131").Value = 940.7143
$ws.Range("J131").Value = 10401.846
$ws.Range("K131").Value = 2822.1429
$ws.Range("L131").Value = 31205.538
$ws.Range("M131").Value = 2217.8571
$ws.Range("N131").Value = -41285.538

$ws.Range("H132").Value = 3356.1304
$ws.Range("I132").Value = 2769.8
$ws.Range("J132").Value = 3807.1538
$ws.Range("K132").Value = 24928.2
$ws.Range("L132").Value = 34264.3842
$ws.Range("M132").Value = -22398.2
$ws.Range("N132").Value = -39324.3842

$ws.Range("H139").Value = 5999
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 5999
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 17997
$ws.Range("N139").Value = -28277
$ws.Range("M139").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1063.75
$ws.Range("I16").Value = 1196.3334
$ws.Range("J16").Value = 666
$ws.Range("K16").Value = 1196.3334
$ws.Range("L16").Value = 666
$ws.Range("M16").Value = -1026.3334
$ws.Range("N16").Value = -1006

$ws.Range("H61").Value = 9956.833000000001
$ws.Range("I61").Value = 12249.956
$ws.Range("J61").Value = 2422.2856
$ws.Range("K61").Value = 12249.956
$ws.Range("L61").Value = 2422.2856
$ws.Range("M61").Value = -12047.956
$ws.Range("N61").Value = -2826.2856

$ws.Range("I68").Value = 5833.3335
$ws.Range("J68").Value = 12500
$ws.Range("K68").Value = 5833.3335
$ws.Range("L68").Value = 12500
$ws.Range("M68").Value = -5084.3335
$ws.Range("N68").Value = -13998

$ws.Range("I71").Value = 5833.3335
$ws.Range("J71").Value = 12500
$ws.Range("K71").Value = 29166.6675
$ws.Range("L71").Value = 62500
$ws.Range("M71").Value = -25422.6675
$ws.Range("N71").Value = -69988

$ws.Range("H113").Value = 9956.833000000001
$ws.Range("I113").Value = 12249.956
$ws.Range("J113").Value = 2422.2856
$ws.Range("K113").Value = 12249.956
$ws.Range("L113").Value = 2422.2856
$ws.Range("M113").Value = -10079.956
$ws.Range("N113").Value = -6762.2856

$ws.Range("H132").Value = 992660.5600000001
$ws.Range("I132").Value = 1333748.2
$ws.Range("J132").Value = 7296
$ws.Range("K132").Value = 4001244.6
$ws.Range("L132").Value = 21888
$ws.Range("M132").Value = -3998714.6
$ws.Range("N132").Value = -26948

$ws.Range("H136").Value = 11066.444
$ws.Range("I136").Value = 9933
$ws.Range("J136").Value = 13333.333
$ws.Range("K136").Value = 29799
$ws.Range("L136").Value = 39999.999
$ws.Range("M136").Value = -27249
$ws.Range("N136").Value = -45099.999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("M62").ClearContents()
$ws.Range("N62").ClearContents()

$ws.Range("H63").Value = 50161.832
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 50161.832
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 50161.832
$ws.Range("N63").Value = -51409.832

$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("M65").ClearContents()
$ws.Range("N65").ClearContents()

$ws.Range("H66").Value = 50161.832
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 50161.832
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 150485.496
$ws.Range("N66").Value = -156725.496

$ws.Range("I132").Value = 2659889.8
$ws.Range("J132").Value = 44550450
$ws.Range("K132").Value = 7979669.399999999
$ws.Range("L132").Value = 133651350
$ws.Range("M132").Value = -7977139.399999999
$ws.Range("N132").Value = -133656410
